$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -0.9383395006683993
$ws.Range("D2").Value = 0.3582583152404712

$ws.Range("C3").Value = 0.1373251200596857
$ws.Range("D3").Value = 0.8920228717674139

$ws.Range("C4").Value = -0.4441432869990896
$ws.Range("D4").Value = 0.6612770685126521

$ws.Range("C5").Value = 0.09165742959291921
$ws.Range("D5").Value = 0.9277997037381802

$ws.Range("C6").Value = 1.117102680396532
$ws.Range("D6").Value = 0.2760023552355115

$ws.Range("C7").Value = 0.815981944732308
$ws.Range("D7").Value = 0.4232578134888247

$ws.Range("C8").Value = 1.422601133002142
$ws.Range("D8").Value = 0.1688771291495026

$ws.Range("C9").Value = -0.5800589973525314
$ws.Range("D9").Value = 0.5677688305398894

$ws.Range("C10").Value = -0.054918835384408
$ws.Range("D10").Value = 0.9566988197265618

$ws.Range("C11").Value = 0.6028083558519614
$ws.Range("D11").Value = 0.552800091285425
